$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Merge the split "Sat Sep 22" / " 11:39:13 PDT 2017" runs into one
#    run by replacing the full matched text with itself.
# ------------------------------------------------------------------
$found = $d.Content.Find.Execute(
    "Sat Sep 22 11:39:13 PDT 2017", $false, $false, $false, $false, $false,
    $true, 1, $false, "Sat Sep 22 11:39:13 PDT 2017", 2)
if (-not $found) {
    throw "Could not find the 'Sat Sep 22 11:39:13 PDT 2017' text to merge"
}

# ------------------------------------------------------------------
# 2) Insert the new "24/09/2017 MAMTHA CHITRA CHICK IN" record block
#    right after the "Amount balance ... - 27842.0" paragraph, and
#    before the first blank paragraph that follows it.
# ------------------------------------------------------------------
$anchor = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*- 27842.0*") {
        $anchor = $p
    }
}
if ($null -eq $anchor) {
    throw "Could not find the '- 27842.0' anchor paragraph"
}

function New-BlankParagraph($afterParagraph, $bold) {
    $afterParagraph.Range.InsertParagraphAfter()
    $newIndex = $afterParagraph.Index + 1
    $newP = $d.Paragraphs.Item($newIndex)
    $newP.Range.Text = "ZPLACEHOLDER"
    $newP.Range.Font.Bold = $bold
    $ok = $d.Content.Find.Execute("ZPLACEHOLDER", $false, $false, $false, $false, $false,
                                   $true, 1, $false, "", 2)
    if (-not $ok) {
        throw "Failed to clear placeholder text for blank paragraph"
    }
    return $newP
}

function New-TextParagraph($afterParagraph, $text, $bold, $color) {
    $afterParagraph.Range.InsertParagraphAfter()
    $newIndex = $afterParagraph.Index + 1
    $newP = $d.Paragraphs.Item($newIndex)
    $newP.Range.Font.Bold = $bold
    $newP.Range.Font.Color = $color
    $newP.Range.Text = $text
    return $newP
}

$autoColor = -16777216
$redColor = 255

# Paragraph: blank, bold
$cur = New-BlankParagraph $anchor $true

# Paragraph: date line
$txt = "Sun Sep 23 10:53:05 PDT 2017"
$cur = New-TextParagraph $cur $txt $false $autoColor

# Paragraph: Person Name ... - HG
$txt = "Person Name" + "`t`t`t`t" + "- HG"
$cur = New-TextParagraph $cur $txt $false $autoColor

# Paragraph: dashed separator
$txt = "---------------------------------------------------------------"
$cur = New-TextParagraph $cur $txt $false $autoColor

# Paragraph: Item Name ... - CARROT
$txt = "Item Name" + "`t`t`t`t" + "- CARROT"
$cur = New-TextParagraph $cur $txt $false $autoColor

# Paragraph: Amount Received ... - 3500 (red)
$txt = "Amount Received" + "`t`t`t" + "- 3500"
$cur = New-TextParagraph $cur $txt $false $redColor

# Paragraph: Amount balance ... - 24342.0 (bold)
$txt = "Amount balance" + "`t`t`t" + "- 24342.0"
$cur = New-TextParagraph $cur $txt $true $autoColor

# Paragraph: Amount Received mode ... - CASH
$txt = "Amount Received mode" + "`t`t" + "- CASH"
$cur = New-TextParagraph $cur $txt $false $autoColor

# Paragraph: blank, not bold
$cur = New-BlankParagraph $cur $false

# Paragraph: blank, bold
$cur = New-BlankParagraph $cur $true

Write-Output "Edit complete"
